$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the "department" header to "departmentNames" (F1) to fix the CSV
# import for department names that contain commas and spaces.
$ws.Range("F1").Value = "departmentNames"

# Restore a sensible active selection on the sheet.
$ws.Range("F1").Select()
